$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the shared string text used by cell C2's lookup ("Result" column value)
$ws.Range("C2").Value = "Sucess! Calendar Invite Sent"

# Remove the stray row 15 (cell C15) that held the old status text
$ws.Rows.Item(15).Delete()

# Reset the active selection to C2
$ws.Range("C2").Select()
